# Adds a new Sheet4 with two small "ListObject" tables describing
# test cases for a Trouble(double a, double b) function, and makes
# it the active sheet.

$wb = $excel.ActiveWorkbook

# Add a new worksheet after Sheet3 and rename it Sheet4
$sheet3 = $wb.Worksheets.Item("Sheet3")
$ws = $wb.Worksheets.Add([System.Type]::Missing, $sheet3)
$ws.Name = "Sheet4"

# ---- First table: C4:G6 ----
$ws.Range("C4").Value = "t1"
$ws.Range("D4").Value = "t2"
$ws.Range("E4").Value = "a"
$ws.Range("F4").Value = "b"
$ws.Range("G4").Value = "Trouble(double a,double b)"

$ws.Range("C5").Value = "Touble(5,5) = 5"
$ws.Range("D5").Value = "Trouble(4,4) = 4"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 6
$ws.Range("G5").Value = "Trouble(5,4)  = 4.5"

$ws.Range("C6").Value = "Touble(7,3) = 5"
$ws.Range("D6").Value = "Trouble(4,4) = 4"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = 4
$ws.Range("G6").Value = 4.5

# ---- Second table: C10:G16 ----
$ws.Range("C10").Value = "t1"
$ws.Range("D10").Value = "t2"
$ws.Range("E10").Value = "a"
$ws.Range("F10").Value = "b"
$ws.Range("G10").Value = "Trouble(double a,double b)"

$ws.Range("C11").Value = "Touble(2,5) = 4"
$ws.Range("D11").Value = "Trouble(1,4) = 4"
$ws.Range("E11").Value = 0
$ws.Range("F11").Value = 6
$ws.Range("G11").Value = "…"

$ws.Range("C12").Value = "Trouble(4,4) = 4"
$ws.Range("D12").Value = "…"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = 5
$ws.Range("G12").Value = "…"

$ws.Range("C13").Value = 4
$ws.Range("D13").Value = "…"
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = 4
$ws.Range("G13").Value = "…"

$ws.Range("C14").Value = 4
$ws.Range("D14").Value = "Trouble(2,2) = 4"
$ws.Range("E14").Value = 1
$ws.Range("F14").Value = 4
$ws.Range("G14").Value = "…"

$ws.Range("C15").Value = 4
$ws.Range("D15").Value = 2
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = "Trouble(4,2) = 3"

$ws.Range("C16").Value = "null"
$ws.Range("D16").Value = "null"
$ws.Range("E16").Value = 4
$ws.Range("F16").Value = 2
$ws.Range("G16").Value = 3

# ---- Column widths to match the other "table" sheets ----
$ws.Columns.Item(1).ColumnWidth = 9.140625
$ws.Columns.Item(2).ColumnWidth = 9.140625
$ws.Columns.Item(3).ColumnWidth = 15.7109375
$ws.Columns.Item(4).ColumnWidth = 15.7109375
$ws.Columns.Item(5).ColumnWidth = 11
$ws.Columns.Item(6).ColumnWidth = 11
$ws.Columns.Item(7).ColumnWidth = 29.7109375
$ws.Columns.Item(8).ColumnWidth = 11
$ws.Columns.Item(9).ColumnWidth = 11

# ---- Turn the two ranges into Excel Tables (ListObjects) ----
$lo1 = $ws.ListObjects.Add(1, $ws.Range("C4:G6"), [System.Type]::Missing, 1)
$lo1.Name = "Table7"

$lo2 = $ws.ListObjects.Add(1, $ws.Range("C10:G16"), [System.Type]::Missing, 1)
$lo2.Name = "Table79"

# Center-align the header rows like the other tables in the workbook
$ws.Range("C4:G4").HorizontalAlignment = -4108
$ws.Range("C10:G10").HorizontalAlignment = -4108

# ---- Sheet view ----
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Range("G21").Select()

# Make Sheet4 the active / selected tab
$ws.Activate()
